# Auto-generated edit script applying cell value updates per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 8400465
$ws.Range("I33").Value = 10983408
$ws.Range("J33").Value = 5899.75
$ws.Range("K33").Value = 10983408
$ws.Range("L33").Value = 5899.75
$ws.Range("M33").Value = -10983179
$ws.Range("N33").Value = -6357.75
$ws.Range("H74").Value = 10771.777
$ws.Range("I74").Value = 11052.471
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 11052.471
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -10116.471
$ws.Range("N74").Value = -7872
$ws.Range("H77").Value = 10771.777
$ws.Range("I77").Value = 11052.471
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 55262.355
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -50582.355
$ws.Range("N77").Value = -39360
$ws.Range("H92").Value = 1340063.2
$ws.Range("I92").Value = 744604.5600000001
$ws.Range("J92").Value = 2233251.2
$ws.Range("K92").Value = 744604.5600000001
$ws.Range("L92").Value = 2233251.2
$ws.Range("M92").Value = -743356.5600000001
$ws.Range("N92").Value = -2235747.2
$ws.Range("H132").Value = 47858.5
$ws.Range("I132").Value = 58182.277
$ws.Range("J132").Value = 1401.5
$ws.Range("K132").Value = 174546.831
$ws.Range("L132").Value = 4204.5
$ws.Range("M132").Value = -172016.831
$ws.Range("N132").Value = -9264.5
$ws.Range("H135").Value = 1304.4
$ws.Range("I135").Value = 1304.4
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 11739.6
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -9204.6
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 9094177
$ws.Range("I137").Value = 1967.75
$ws.Range("K137").Value = 5903.25
$ws.Range("M137").Value = -3353.25
$ws.Range("H138").Value = 6550.418
$ws.Range("I138").Value = 10638.389
$ws.Range("J138").Value = 4561.676
$ws.Range("K138").Value = 31915.167
$ws.Range("L138").Value = 13685.028
$ws.Range("M138").Value = -26775.167
$ws.Range("N138").Value = -23965.028

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 107607.33
$ws.Range("I32").Value = 111989.79
$ws.Range("J32").Value = 9002
$ws.Range("K32").Value = 111989.79
$ws.Range("L32").Value = 9002
$ws.Range("M32").Value = -111702.79
$ws.Range("N32").Value = -9576
$ws.Range("H45").Value = 203199.8
$ws.Range("I45").Value = 203199.8
$ws.Range("K45").Value = 203199.8
$ws.Range("M45").Value = -202822.8
$ws.Range("H61").Value = 1073749.2
$ws.Range("I61").Value = 3637.4443
$ws.Range("J61").Value = 5889252.5
$ws.Range("K61").Value = 3637.4443
$ws.Range("L61").Value = 5889252.5
$ws.Range("M61").Value = -3425.4443
$ws.Range("N61").Value = -5889676.5
$ws.Range("H74").Value = 971853.2
$ws.Range("I74").Value = 2099.4
$ws.Range("J74").Value = 1241229.2
$ws.Range("K74").Value = 2099.4
$ws.Range("L74").Value = 1241229.2
$ws.Range("M74").Value = -1225.4
$ws.Range("N74").Value = -1242977.2
$ws.Range("H77").Value = 971853.2
$ws.Range("I77").Value = 2099.4
$ws.Range("J77").Value = 1241229.2
$ws.Range("K77").Value = 10497
$ws.Range("L77").Value = 6206146
$ws.Range("M77").Value = -6129
$ws.Range("N77").Value = -6214882
$ws.Range("H97").Value = 11828
$ws.Range("I97").Value = 12975.556
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 12975.556
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -12479.556
$ws.Range("N97").Value = -2492
$ws.Range("H102").Value = 1855.619
$ws.Range("I102").Value = 1798.4
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1798.4
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -176.4000000000001
$ws.Range("N102").Value = -6244
$ws.Range("H132").Value = 5633.3335
$ws.Range("I132").Value = 4900
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 14700
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -12170
$ws.Range("N132").Value = -23060
$ws.Range("H136").Value = 1073749.2
$ws.Range("I136").Value = 3637.4443
$ws.Range("J136").Value = 5889252.5
$ws.Range("K136").Value = 10912.3329
$ws.Range("L136").Value = 17667757.5
$ws.Range("M136").Value = -8362.332900000001
$ws.Range("N136").Value = -17672857.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1087.6666
$ws.Range("I20").Value = 990.875
$ws.Range("J20").Value = 1147.2307
$ws.Range("K20").Value = 990.875
$ws.Range("L20").Value = 1147.2307
$ws.Range("M20").Value = -743.875
$ws.Range("N20").Value = -1641.2307
$ws.Range("H94").Value = 2781466
$ws.Range("I94").Value = 3836223.5
$ws.Range("J94").Value = 741.7273
$ws.Range("K94").Value = 3836223.5
$ws.Range("L94").Value = 741.7273
$ws.Range("M94").Value = -3835772.5
$ws.Range("N94").Value = -1643.7273
$ws.Range("H99").Value = 8060.8096
$ws.Range("I99").Value = 14057.8
$ws.Range("K99").Value = 14057.8
$ws.Range("M99").Value = -12559.8
$ws.Range("H134").Value = 18751814
$ws.Range("I134").Value = 1119.1842
$ws.Range("J134").Value = 90004450
$ws.Range("K134").Value = 3357.5526
$ws.Range("L134").Value = 270013350
$ws.Range("M134").Value = -822.5526
$ws.Range("N134").Value = -270018420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3831.5789
$ws.Range("I16").Value = 2510.6365
$ws.Range("J16").Value = 5647.875
$ws.Range("K16").Value = 2510.6365
$ws.Range("L16").Value = 5647.875
$ws.Range("M16").Value = -2223.6365
$ws.Range("N16").Value = -6221.875
$ws.Range("H22").Value = 605.21277
$ws.Range("I22").Value = 462.07318
$ws.Range("J22").Value = 1583.3334
$ws.Range("K22").Value = 462.07318
$ws.Range("L22").Value = 1583.3334
$ws.Range("M22").Value = -112.07318
$ws.Range("N22").Value = -2283.3334
$ws.Range("H31").Value = 5237.269
$ws.Range("I31").Value = 3210.5
$ws.Range("J31").Value = 6138.0557
$ws.Range("K31").Value = 3210.5
$ws.Range("L31").Value = 6138.0557
$ws.Range("M31").Value = -2915.5
$ws.Range("N31").Value = -6728.0557
$ws.Range("H34").Value = 5237.269
$ws.Range("I34").Value = 3210.5
$ws.Range("J34").Value = 6138.0557
$ws.Range("K34").Value = 3210.5
$ws.Range("L34").Value = 6138.0557
$ws.Range("M34").Value = -3008.5
$ws.Range("N34").Value = -6542.0557
$ws.Range("H113").Value = 3831.5789
$ws.Range("I113").Value = 2510.6365
$ws.Range("J113").Value = 5647.875
$ws.Range("K113").Value = 2510.6365
$ws.Range("L113").Value = 5647.875
$ws.Range("M113").Value = -340.6365000000001
$ws.Range("N113").Value = -9987.875
$ws.Range("H122").Value = 2733.7222
$ws.Range("I122").Value = 2948.9
$ws.Range("J122").Value = 2464.75
$ws.Range("K122").Value = 8846.700000000001
$ws.Range("L122").Value = 7394.25
$ws.Range("M122").Value = -6396.700000000001
$ws.Range("N122").Value = -12294.25
$ws.Range("H132").Value = 3059.5
$ws.Range("I132").Value = 2881.8064
$ws.Range("J132").Value = 3671.5557
$ws.Range("K132").Value = 8645.4192
$ws.Range("L132").Value = 11014.6671
$ws.Range("M132").Value = -6115.4192
$ws.Range("N132").Value = -16074.6671
$ws.Range("H134").Value = 1983.6364
$ws.Range("I134").Value = 1379
$ws.Range("J134").Value = 4229.4287
$ws.Range("K134").Value = 4137
$ws.Range("L134").Value = 12688.2861
$ws.Range("M134").Value = -1602
$ws.Range("N134").Value = -17758.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 4454352.5
$ws.Range("I11").Value = 82325.55
$ws.Range("K11").Value = 246976.65
$ws.Range("M11").Value = -246836.65
$ws.Range("H131").Value = 2485723.5
$ws.Range("I131").Value = 4133239.5
$ws.Range("J131").Value = 69366.47
$ws.Range("K131").Value = 12399718.5
$ws.Range("L131").Value = 208099.41
$ws.Range("M131").Value = -12394678.5
$ws.Range("N131").Value = -218179.41

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8185279.5
$ws.Range("I80").Value = 158320.16
$ws.Range("J80").Value = 20065180
$ws.Range("K80").Value = 158320.16
$ws.Range("L80").Value = 20065180
$ws.Range("M80").Value = -157322.16
$ws.Range("N80").Value = -20067176
$ws.Range("H83").Value = 8185279.5
$ws.Range("I83").Value = 158320.16
$ws.Range("J83").Value = 20065180
$ws.Range("K83").Value = 791600.8
$ws.Range("L83").Value = 100325900
$ws.Range("M83").Value = -786608.8
$ws.Range("N83").Value = -100335884
$ws.Range("H113").Value = 1727.8125
$ws.Range("I113").Value = 1709.6666
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1709.6666
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 460.3334
$ws.Range("N113").Value = -6340
$ws.Range("H132").Value = 613782
$ws.Range("I132").Value = 2881.0356
$ws.Range("J132").Value = 1835583.9
$ws.Range("K132").Value = 8643.106800000001
$ws.Range("L132").Value = 5506751.699999999
$ws.Range("M132").Value = -6113.106800000001
$ws.Range("N132").Value = -5511811.699999999
$ws.Range("H139").Value = 500000
$ws.Range("J139").Value = 500000
$ws.Range("L139").Value = 500000
$ws.Range("N139").Value = -510280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1391353.6
$ws.Range("I40").Value = 1636008.1
$ws.Range("J40").Value = 4978.1665
$ws.Range("K40").Value = 1636008.1
$ws.Range("L40").Value = 4978.1665
$ws.Range("M40").Value = -1635872.1
$ws.Range("N40").Value = -5250.1665
$ws.Range("H61").Value = 5002221.5
$ws.Range("J61").Value = 2061.923
$ws.Range("L61").Value = 2061.923
$ws.Range("N61").Value = -2465.923
$ws.Range("H93").Value = 3024.7144
$ws.Range("I93").Value = 2776.9092
$ws.Range("J93").Value = 3933.3333
$ws.Range("K93").Value = 2776.9092
$ws.Range("L93").Value = 3933.3333
$ws.Range("M93").Value = -1528.9092
$ws.Range("N93").Value = -6429.3333
$ws.Range("H113").Value = 5002221.5
$ws.Range("J113").Value = 2061.923
$ws.Range("L113").Value = 2061.923
$ws.Range("N113").Value = -6401.923
$ws.Range("H127").Value = 133326.67
$ws.Range("J127").Value = 133326.67
$ws.Range("L127").Value = 133326.67
$ws.Range("N127").Value = -143246.67
$ws.Range("H132").Value = 28997.75
$ws.Range("I132").Value = 2991
$ws.Range("J132").Value = 37666.668
$ws.Range("K132").Value = 8973
$ws.Range("L132").Value = 113000.004
$ws.Range("M132").Value = -6443
$ws.Range("N132").Value = -118060.004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1567.4546
$ws.Range("I122").Value = 1567.4546
$ws.Range("K122").Value = 4702.3638
$ws.Range("M122").Value = -2252.3638
$ws.Range("H132").Value = 103874.9
$ws.Range("I132").Value = 501224.5
$ws.Range("J132").Value = 4537.5
$ws.Range("K132").Value = 1503673.5
$ws.Range("L132").Value = 13612.5
$ws.Range("M132").Value = -18672.5

